$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.349.82"
$ws.Range("E2").Value = "  -0.60%  "

$ws.Range("D3").Value = "1.568.54"
$ws.Range("E3").Value = "  +0.29%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.39%  "

$ws.Range("E6").Value = "  -0.62%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.23"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.81"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.48%  "

$ws.Range("E10").Value = "  -0.92%  "

$ws.Range("E11").Value = "  -0.84%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0894"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.28%  "

$ws.Range("D13").Value = "1.793.77"
$ws.Range("E13").Value = "  +0.42%  "

$ws.Range("D14").Value = "1.566.70"
$ws.Range("E14").Value = "  +1.55%  "

$ws.Range("E15").Value = "  -0.41%  "

$ws.Range("D16").Value = "28.332.99"
$ws.Range("E16").Value = "  -0.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.513"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.52%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "227.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.07%  "

$ws.Range("E20").Value = "  +0.55%  "

$ws.Range("E21").Value = "  -1.94%  "

$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.96%  "

$ws.Range("E25").Value = "  -0.60%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.49%  "

$ws.Range("E28").Value = "  -0.58%  "

$ws.Range("E29").Value = "  -1.59%  "

$ws.Range("E30").Value = "  +0.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0479"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.84%  "

$ws.Range("E33").Value = "  -0.87%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.17%  "

$ws.Range("D35").Value = "1.378.49"
$ws.Range("E35").Value = "  -1.20%  "

$ws.Range("E36").Value = "  +2.12%  "

$ws.Range("E37").Value = "  -3.12%  "

$ws.Range("E38").Value = "  -0.23%  "

$ws.Range("E39").Value = "  +2.54%  "

$ws.Range("E40").Value = "  -2.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.519"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.92"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.43%  "

$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0475"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("E45").Value = "  -0.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "62.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.00%  "

$ws.Range("E48").Value = "  -6.27%  "

$ws.Range("D49").Value = "1.705.17"
$ws.Range("E49").Value = "  +0.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "85.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.91%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0101"
$ws.Range("E51").Value = "  -0.92%  "
